# Adding randomization of test run, tested.
# Append three more rows of generated test data to the "userId" sheet
# (rows 6-8: userId 5,6,7 all with #enable# = TRUE), then make that
# sheet the active/selected sheet with B6:B8 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userId")

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = $true

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = $true

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = $true

# Make the userId sheet the active sheet and select the newly added
# boolean values, matching the updated selection/tabSelected state.
$ws.Activate()
$ws.Range("B6:B8").Select()
